# Fruta / hortaliza, semanal
# Insert a new weekly price-report pair of rows (Primera / Segunda quality)
# for "Perejil" at "Vega Monumental Concepción" right before the existing
# block that starts at row 224, shifting the rest of the table down by two
# rows (224-237 -> 226-239).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 224-225; everything from the old row 224 onward
# (through 237) shifts down to 226-239.
$ws.Range("A224:A225").EntireRow.Insert()

# New row 224 - Calidad "Primera"
$ws.Cells.Item(224, 1).Value = 11
$ws.Cells.Item(224, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(224, 3).Value = "Bíobío"
$ws.Cells.Item(224, 4).Value = 45147
$ws.Cells.Item(224, 5).Value = 8
$ws.Cells.Item(224, 6).Value = 100112044
$ws.Cells.Item(224, 7).Value = "Perejil"
$ws.Cells.Item(224, 8).Value = "Sin especificar"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 200
$ws.Cells.Item(224, 11).Value = 600
$ws.Cells.Item(224, 12).Value = 700
$ws.Cells.Item(224, 13).Value = 650
$ws.Cells.Item(224, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(224, 15).Value = "Región de Ñuble"
$ws.Cells.Item(224, 16).Value = 650
$ws.Cells.Item(224, 17).Value = 1
$ws.Cells.Item(224, 18).Value = "Hortaliza"

# New row 225 - Calidad "Segunda"
$ws.Cells.Item(225, 1).Value = 11
$ws.Cells.Item(225, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(225, 3).Value = "Bíobío"
$ws.Cells.Item(225, 4).Value = 45147
$ws.Cells.Item(225, 5).Value = 8
$ws.Cells.Item(225, 6).Value = 100112044
$ws.Cells.Item(225, 7).Value = "Perejil"
$ws.Cells.Item(225, 8).Value = "Sin especificar"
$ws.Cells.Item(225, 9).Value = "Segunda"
$ws.Cells.Item(225, 10).Value = 100
$ws.Cells.Item(225, 11).Value = 500
$ws.Cells.Item(225, 12).Value = 500
$ws.Cells.Item(225, 13).Value = 500
$ws.Cells.Item(225, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(225, 15).Value = "Región de Ñuble"
$ws.Cells.Item(225, 16).Value = 500
$ws.Cells.Item(225, 17).Value = 1
$ws.Cells.Item(225, 18).Value = "Hortaliza"
